$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row containing "「声を上げろ」..." (row 209) was removed from the
# data, shifting all subsequent rows up by one.
$ws.Rows.Item(209).Delete()
